$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at D:E, shifting flexibility..class right to F..L
$ws.Range("D:E").Insert()

# New header labels for the inserted sum_SASA / max_SASA columns
$ws.Range("D1").Value = "sum_SASA"
$ws.Range("E1").Value = "max_SASA"

# Populate sum_SASA (D) / max_SASA (E) data rows
$ws.Range("D2").Value = 3.938433140488909
$ws.Range("E2").Value = 3.938433140488909
$ws.Range("D3").Value = 3.838393281831593
$ws.Range("E3").Value = 3.838393281831593
$ws.Range("D4").Value = 4.037903583503363
$ws.Range("E4").Value = 4.037903583503363
$ws.Range("D5").Value = 3.941532126214653
$ws.Range("E5").Value = 3.941532126214653
$ws.Range("D6").Value = 3.850616806854433
$ws.Range("E6").Value = 3.850616806854433
$ws.Range("D7").Value = 4.157103523039212
$ws.Range("E7").Value = 4.157103523039212
$ws.Range("D8").Value = 3.795369541055058
$ws.Range("E8").Value = 3.795369541055058
$ws.Range("D9").Value = 3.815165399585339
$ws.Range("E9").Value = 3.815165399585339
$ws.Range("D10").Value = 4.148720083623762
$ws.Range("E10").Value = 4.148720083623762
$ws.Range("D11").Value = 4.267455298719907
$ws.Range("E11").Value = 4.267455298719907
$ws.Range("D12").Value = 4.021427922899613
$ws.Range("E12").Value = 4.021427922899613
$ws.Range("D13").Value = 4.109653098323708
$ws.Range("E13").Value = 4.109653098323708
$ws.Range("D14").Value = 3.432590204636147
$ws.Range("E14").Value = 3.432590204636147
$ws.Range("D15").Value = 7.655021391887416
$ws.Range("E15").Value = 3.918236369513348
$ws.Range("D16").Value = 7.178448819228891
$ws.Range("E16").Value = 3.952893427540636
$ws.Range("D17").Value = 3.982816231500002
$ws.Range("E17").Value = 3.982816231500002
$ws.Range("D18").Value = 4.379913696596782
$ws.Range("E18").Value = 4.379913696596782
$ws.Range("D19").Value = 4.536346695756726
$ws.Range("E19").Value = 4.536346695756726
$ws.Range("D20").Value = 3.829348253889487
$ws.Range("E20").Value = 3.829348253889487
$ws.Range("D21").Value = 8.005216139406006
$ws.Range("E21").Value = 4.097181451668026
